$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Directors model: add YearBorn property (reuses existing shared string)
$ws.Range("B3").Value = "YearBorn"

# Studios model: add Trivia property (new shared string, must be created
# before "Short Bio" so shared-string indices match the target ordering)
$ws.Range("C5").Value = "Trivia"

# Directors model: add Short Bio property (new shared string)
$ws.Range("D4").Value = "Short Bio"

# Move active selection to B4
$ws.Range("B4").Select()
